# GPLIM-5135 fixes and improvements for Pooled Tube Upload
#
# The sample row (row 2) on Sheet1 is missing a value in the "Volume"
# column (P). Populate it, and leave the active selection on that cell,
# matching the state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 61

[void]$ws.Range("P2").Select()
